$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '35.485.12'
$ws.Cells.Item(2, 5).Value = '  +1.17%  '
$ws.Cells.Item(3, 4).Value = '1.902.88'
$ws.Cells.Item(3, 5).Value = '  +2.65%  '
$ws.Cells.Item(4, 5).Value = '  +0.53%  '
$ws.Cells.Item(5, 4).Value = "'245.57"
$ws.Cells.Item(5, 5).Value = '  +4.23%  '
$ws.Cells.Item(6, 5).Value = '  +1.46%  '
$ws.Cells.Item(7, 5).Value = '  +0.43%  '
$ws.Cells.Item(8, 4).Value = "'41.89"
$ws.Cells.Item(8, 5).Value = '  -2.35%  '
$ws.Cells.Item(9, 4).Value = "'0.340"
$ws.Cells.Item(9, 5).Value = '  +2.94%  '
$ws.Cells.Item(10, 4).Value = "'0.0703"
$ws.Cells.Item(10, 5).Value = '  +1.15%  '
$ws.Cells.Item(11, 5).Value = '  +1.09%  '
$ws.Cells.Item(12, 4).Value = '2.177.80'
$ws.Cells.Item(12, 5).Value = '  +2.64%  '
$ws.Cells.Item(13, 4).Value = "'12.34"
$ws.Cells.Item(13, 5).Value = '  +8.20%  '
$ws.Cells.Item(14, 2).Value = 'Polygon'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Cells.Item(14, 4).Value = "'0.690"
$ws.Cells.Item(14, 5).Value = '  +1.28%  '
$ws.Cells.Item(15, 2).Value = 'Polkadot'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Cells.Item(15, 4).Value = "'4.86"
$ws.Cells.Item(15, 5).Value = '  +3.54%  '
$ws.Cells.Item(16, 2).Value = 'WrappedEther'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(16, 4).Value = '1.884.20'
$ws.Cells.Item(16, 5).Value = '  +1.53%  '
$ws.Cells.Item(17, 4).Value = '35.492.59'
$ws.Cells.Item(17, 5).Value = '  +1.19%  '
$ws.Cells.Item(18, 4).Value = "'71.85"
$ws.Cells.Item(18, 5).Value = '  +2.27%  '
$ws.Cells.Item(19, 4).Value = '0.0₃0822'
$ws.Cells.Item(19, 5).Value = '  +3.10%  '
$ws.Cells.Item(20, 4).Value = "'242.95"
$ws.Cells.Item(20, 5).Value = '  +0.63%  '
$ws.Cells.Item(21, 4).Value = "'12.50"
$ws.Cells.Item(22, 5).Value = '  +1.67%  '
$ws.Cells.Item(23, 5).Value = '  +0.45%  '
$ws.Cells.Item(24, 4).Value = "'2.30"
$ws.Cells.Item(24, 5).Value = '  +1.33%  '
$ws.Cells.Item(25, 4).Value = "'172.22"
$ws.Cells.Item(25, 5).Value = '  +0.36%  '
$ws.Cells.Item(26, 4).Value = "'2.18"
$ws.Cells.Item(26, 5).Value = '  +17.40%  '
$ws.Cells.Item(27, 4).Value = "'8.55"
$ws.Cells.Item(27, 5).Value = '  +7.92%  '
$ws.Cells.Item(28, 4).Value = "'17.95"
$ws.Cells.Item(28, 5).Value = '  +1.34%  '
$ws.Cells.Item(29, 5).Value = '  +0.24%  '
$ws.Cells.Item(30, 4).Value = "'0.970"
$ws.Cells.Item(30, 5).Value = '  +23.14%  '
$ws.Cells.Item(31, 4).Value = "'0.0569"
$ws.Cells.Item(31, 5).Value = '  +2.06%  '
$ws.Cells.Item(32, 4).Value = "'4.10"
$ws.Cells.Item(32, 5).Value = '  +2.34%  '
$ws.Cells.Item(33, 5).Value = '  +0.49%  '
$ws.Cells.Item(34, 4).Value = "'4.16"
$ws.Cells.Item(34, 5).Value = '  +4.40%  '
$ws.Cells.Item(35, 5).Value = '  +7.50%  '
$ws.Cells.Item(36, 5).Value = '  +10.75%  '
$ws.Cells.Item(37, 4).Value = "'2.02"
$ws.Cells.Item(37, 5).Value = '  -2.19%  '
$ws.Cells.Item(38, 5).Value = '  +2.62%  '
$ws.Cells.Item(39, 5).Value = '  +1.42%  '
$ws.Cells.Item(40, 4).Value = "'90.78"
$ws.Cells.Item(40, 5).Value = '  -1.26%  '
$ws.Cells.Item(41, 4).Value = "'15.73"
$ws.Cells.Item(41, 5).Value = '  +4.71%  '
$ws.Cells.Item(42, 4).Value = '1.348.97'
$ws.Cells.Item(42, 5).Value = '  -0.10%  '
$ws.Cells.Item(43, 4).Value = "'0.0616"
$ws.Cells.Item(43, 5).Value = '  +14.11%  '
$ws.Cells.Item(44, 4).Value = "'49.42"
$ws.Cells.Item(44, 5).Value = '  +42.16%  '
$ws.Cells.Item(45, 2).Value = 'RenderToken'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(45, 4).Value = "'2.35"
$ws.Cells.Item(45, 5).Value = '  +1.24%  '
$ws.Cells.Item(46, 2).Value = 'Gas'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Cells.Item(46, 4).Value = "'12.95"
$ws.Cells.Item(46, 5).Value = '  +2.46%  '
$ws.Cells.Item(47, 4).Value = "'2.41"
$ws.Cells.Item(47, 5).Value = '  +0.74%  '
$ws.Cells.Item(48, 4).Value = "'2.75"
$ws.Cells.Item(48, 5).Value = '  -0.37%  '
$ws.Cells.Item(49, 5).Value = '  +3.40%  '
$ws.Cells.Item(50, 4).Value = '2.087.07'
$ws.Cells.Item(50, 5).Value = '  +2.56%  '
$ws.Cells.Item(51, 4).Value = "'0.0690"
$ws.Cells.Item(51, 5).Value = '  +1.48%  '
